# Rename the inline picture shapes that live in the document's headers and
# footers, per the authored diff:
#   - the Pearson Edexcel logo picture (alt text / description =
#     "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png")
#     is renamed from "image2.png" to "image1.png"
#   - the BTEC orange logo picture (alt text / description = "BTec_Logo-Orange")
#     is renamed from "image1.jpg" to "image2.jpg"
# These pictures are inline shapes anchored inside the header/footer stories,
# so we walk every section's Headers and Footers collections (both the
# "first page" and "default" header/footer parts) and rename every inline
# shape we find there based on its current description, rather than
# hard-coding which physical header/footer part it happens to live in.

$d = $word.ActiveDocument

function Rename-LogoInlineShapes($range) {
    $ishapes = $range.InlineShapes
    for ($i = 1; $i -le $ishapes.Count; $i++) {
        $shp = $ishapes.Item($i)
        $desc = $shp.AlternativeText

        if ($desc -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
            $shp.Name = "image1.png"
        }
        elseif ($desc -eq "BTec_Logo-Orange") {
            $shp.Name = "image2.jpg"
        }
    }
}

$sections = $d.Sections
for ($si = 1; $si -le $sections.Count; $si++) {
    $sec = $sections.Item($si)

    $headers = $sec.Headers
    for ($hi = 1; $hi -le $headers.Count; $hi++) {
        $h = $headers.Item($hi)
        if ($h.Exists) {
            Rename-LogoInlineShapes $h.Range
        }
    }

    $footers = $sec.Footers
    for ($fi = 1; $fi -le $footers.Count; $fi++) {
        $f = $footers.Item($fi)
        if ($f.Exists) {
            Rename-LogoInlineShapes $f.Range
        }
    }
}
